# GMS Data Release 1
# Refresh the sequencing_report "Name | Type" field listing:
#  - drop the old synthetic "Key" row
#  - rename several fields to their current GMS column names
#  - append a new "data_format" field row
#  - tidy the view (selection + column width) for the new, longer names

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was the placeholder "Key | integer" row - remove it entirely.
# Everything below shifts up by one row after this.
$ws.Rows(2).Delete()

# Row 2 is now participant_id/Integer - unchanged, leave as-is.

# Row 3 was lab_sample_id/Integer -> platekey/varchar
$ws.Range("A3").Value = "platekey"
$ws.Range("B3").Value = "varchar"

# Row 4 was plate_key/varchar -> referral_id/varchar
$ws.Range("A4").Value = "referral_id"

# Row 5 was type/varchar -> associated_interpretation_request_id/varchar
$ws.Range("A5").Value = "associated_interpretation_request_id"

# Row 7 was delivery_date/varchar -> delivery_id/varchar
$ws.Range("A7").Value = "delivery_id"
$ws.Range("B7").Value = "varchar"

# Row 8 was path/varchar -> delivery_date/timestamp
$ws.Range("A8").Value = "delivery_date"
$ws.Range("B8").Value = "timestamp"

# Row 9 was delivery_version/varchar -> path/varchar
$ws.Range("A9").Value = "path"
$ws.Range("B9").Value = "varchar"

# Row 10 was genome_build/varchar -> delivery_version/varchar
$ws.Range("A10").Value = "delivery_version"

# Row 11 no longer exists (sheet only has 10 rows after the deletion) - recreate
# it with the same look as the row above, then fill in genome_build/varchar.
$ws.Range("A10:C10").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$ws.Range("A11").Value = "genome_build"
$ws.Range("B11").Value = "varchar"

# Append the new row 12: data_format/varchar. Both cells share the bordered
# "A-column" look (no yellow highlight column here, so no C12 cell either).
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("A12").Value = "data_format"
$ws.Range("B12").Value = "varchar"

# Row 6 was delivery_id/varchar -> delivery_type/varchar
$ws.Range("A6").Value = "delivery_type"

# Widen column A so the longer field names (e.g.
# "associated_interpretation_request_id") are readable, and leave the
# selection where the editor last left it.
$ws.Columns(1).ColumnWidth = 31.5
$ws.Range("F5").Select() | Out-Null

Write-Output "sequencing_report field list refreshed"
